$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'96.643.31"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = "'3.570.77"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'241.34"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.04%  '
$ws.Range('D6').Value = "'653.08"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('D7').Value = "'1.67"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +14.66%  '
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range('D9').Value = "'1.08"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.73%  '
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').Value = "'3.565.88"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.40%  '
$ws.Range('D12').Value = "'43.59"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.76%  '
$ws.Range('D13').Value = "'0.203"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').Value = "'6.41"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.58%  '
$ws.Range('D15').Value = "'4.234.52"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.62%  '
$ws.Range('D16').Value = "'96.457.03"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.67%  '
$ws.Range('D17').Value = "'0.0000259"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.36%  '
$ws.Range('D18').Value = "'3.580.67"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('D20').Value = "'12.62"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('D21').Value = "'17.99"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').Value = "'0.543"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +13.32%  '
$ws.Range('D23').Value = "'507.95"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('D24').Value = "'3.39"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.89%  '
$ws.Range('D25').Value = "'6.97"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.78%  '
$ws.Range('E26').Value = '  +1.94%  '
$ws.Range('D27').Value = "'96.71"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.41%  '
$ws.Range('D28').Value = "'13.00"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.41%  '
$ws.Range('D29').Value = "'3.761.49"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('D30').Value = "'0.155"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +12.46%  '
$ws.Range('D31').Value = "'3.03"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.49%  '
$ws.Range('D32').Value = "'11.50"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.19%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  +4.39%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = "'31.37"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.37%  '
$ws.Range('D37').Value = "'630.96"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.38%  '
$ws.Range('D38').Value = "'8.88"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.74%  '
$ws.Range('D39').Value = "'0.569"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.23%  '
$ws.Range('E40').Value = '  +11.25%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').Value = "'0.151"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.13%  '
$ws.Range('D43').Value = "'0.905"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.07%  '
$ws.Range('E44').Value = '  +5.63%  '
$ws.Range('D45').Value = "'5.77"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.93%  '
$ws.Range('D46').Value = "'0.0429"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.06%  '
$ws.Range('E47').Value = '  +3.20%  '
$ws.Range('D48').Value = "'23.53"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.66%  '
$ws.Range('D49').Value = "'32.95"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.37%  '
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('D51').Value = "'8.34"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.68%  '
